$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "30.492.51"
$ws.Cells.Item(2, 5).Value = "  -1.63%  "

$ws.Cells.Item(3, 4).Value = "2.098.90"
$ws.Cells.Item(3, 5).Value = "  -0.87%  "

$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  -0.12%  "

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "329.97"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -1.52%  "

$ws.Cells.Item(6, 5).Value = "  -0.03%  "

$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = "0.5224"
$c.Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  -2.20%  "

$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = "0.4430"
$c.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +0.68%  "

$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = "53.60"
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +15.99%  "

$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = "0.08936"
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -1.43%  "

$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "1.153"
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  -2.49%  "

$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = "24.42"
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  -3.55%  "

$ws.Cells.Item(13, 4).Value = "2.091.51"

$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "6.696"
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -1.38%  "

$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "7.718"
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  -1.19%  "

$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "96.36"
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  -1.46%  "

$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = "1.003"
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  +0.06%  "

$ws.Cells.Item(18, 5).Value = "  -1.39%  "

$ws.Cells.Item(19, 5).Value = "  -0.98%  "

$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = "19.15"
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -0.41%  "

$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +0.00%  "

$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "6.280"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -1.74%  "

$ws.Cells.Item(23, 4).Value = "30.529.71"
$ws.Cells.Item(23, 5).Value = "  -1.71%  "

$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "12.29"
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +0.53%  "

$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "2.321"
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +1.84%  "

$ws.Cells.Item(26, 4).Value = "2.335.38"
$ws.Cells.Item(26, 5).Value = "  -1.37%  "

$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "22.28"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -2.65%  "

$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = "2.570"
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  -0.33%  "

$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "163.76"
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +0.00%  "

$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = "131.84"
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -1.81%  "

$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "1.188"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +0.79%  "

$ws.Cells.Item(32, 5).Value = "  -0.73%  "

$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = "1.662"
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +8.50%  "

$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "6.167"
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  -1.55%  "

$ws.Cells.Item(35, 5).Value = "  -2.51%  "

$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = "10.18"
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +6.53%  "

$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = "0.02562"
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -2.38%  "

$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "0.06793"
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +0.43%  "

$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "5.471"
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -1.82%  "

$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "12.72"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -3.13%  "

$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "0.2258"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -1.09%  "

$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "0.6904"
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -0.13%  "

$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = "1.254"
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -0.19%  "

$ws.Cells.Item(44, 5).Value = "  +0.02%  "

$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "14.03"
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -1.18%  "

$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "0.6339"
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -2.97%  "

$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "2.208"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -2.41%  "

$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "3.630"
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -1.50%  "

$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "1.247"
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +6.59%  "

$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "1.246"
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -2.81%  "

$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "81.81"
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -1.75%  "
